# Update the cryptos list (Price / Volume(1h) columns, and the
# WEMIXToken / LidoDAOToken row swap) per the Fri Nov 10 20:52:01 UTC 2023
# GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.276.99'
$ws.Range('E2').Value = '  +2.27%  '
$ws.Range('D3').Value = '2.088.96'
$ws.Range('E3').Value = '  +2.77%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').Value = '251.04'
$ws.Range('E5').Value = '  +2.44%  '
$ws.Range('D6').Value = '0.663'
$ws.Range('E6').Value = '  +0.83%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').Value = '54.75'
$ws.Range('E8').Value = '  +23.80%  '
$ws.Range('D9').Value = '61.67'
$ws.Range('E9').Value = '  +2.71%  '
$ws.Range('D10').Value = '0.379'
$ws.Range('E10').Value = '  +6.23%  '
$ws.Range('D11').Value = '0.0749'
$ws.Range('E11').Value = '  +4.99%  '
$ws.Range('E12').Value = '  +8.05%  '
$ws.Range('D13').Value = '15.03'
$ws.Range('E13').Value = '  +5.51%  '
$ws.Range('D14').Value = '2.391.12'
$ws.Range('E14').Value = '  +3.06%  '
$ws.Range('D15').Value = '0.832'
$ws.Range('E15').Value = '  +3.68%  '
$ws.Range('D16').Value = '2.088.01'
$ws.Range('E16').Value = '  +3.32%  '
$ws.Range('D17').Value = '5.19'
$ws.Range('E17').Value = '  +6.87%  '
$ws.Range('D18').Value = '37.204.12'
$ws.Range('E18').Value = '  +2.46%  '
$ws.Range('D19').Value = '72.84'
$ws.Range('E19').Value = '  +2.95%  '
$ws.Range('D20').Value = '14.47'
$ws.Range('E20').Value = '  +15.92%  '
$ws.Range('D21').Value = '0.0₃0848'
$ws.Range('E21').Value = '  +5.10%  '
$ws.Range('D22').Value = '240.63'
$ws.Range('E22').Value = '  +2.25%  '
$ws.Range('D23').Value = '5.22'
$ws.Range('E23').Value = '  +7.34%  '
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('D25').Value = '2.47'
$ws.Range('E25').Value = '  +1.43%  '
$ws.Range('D26').Value = '172.45'
$ws.Range('E26').Value = '  +2.58%  '
$ws.Range('D27').Value = '9.22'
$ws.Range('E27').Value = '  +6.29%  '
$ws.Range('D28').Value = '20.77'
$ws.Range('E28').Value = '  +4.63%  '
$ws.Range('D29').Value = '2.03'
$ws.Range('E29').Value = '  +5.48%  '
$ws.Range('D30').Value = '0.124'
$ws.Range('E30').Value = '  +2.19%  '
$ws.Range('D31').Value = "'1.10"
$ws.Range('E31').Value = '  +30.59%  '
$ws.Range('D32').Value = '22.53'
$ws.Range('E32').Value = '  +4.37%  '
$ws.Range('D33').Value = '4.53'
$ws.Range('E33').Value = '  +5.02%  '
$ws.Range('D34').Value = "'0.0620"
$ws.Range('E34').Value = '  +7.89%  '
$ws.Range('D35').Value = '0.0891'
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('D36').Value = '0.999'
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('D37').Value = '4.18'
$ws.Range('E37').Value = '  +6.05%  '
$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D38').Value = '2.23'
$ws.Range('E38').Value = '  +2.16%  '
$ws.Range('B39').Value = 'WEMIXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D39').Value = '1.81'
$ws.Range('E39').Value = '  -2.96%  '
$ws.Range('D40').Value = '5.38'
$ws.Range('E40').Value = '  +173.41%  '
$ws.Range('D41').Value = '1.36'
$ws.Range('E41').Value = '  +3.27%  '
$ws.Range('D42').Value = '18.16'
$ws.Range('E42').Value = '  +16.13%  '
$ws.Range('D43').Value = '0.0227'
$ws.Range('E43').Value = '  +7.28%  '
$ws.Range('E44').Value = '  +6.59%  '
$ws.Range('D45').Value = '98.74'
$ws.Range('E45').Value = '  +3.43%  '
$ws.Range('D46').Value = '0.0951'
$ws.Range('E46').Value = '  +17.05%  '
$ws.Range('E47').Value = '  +0.58%  '
$ws.Range('D48').Value = '1.325.19'
$ws.Range('E48').Value = '  +1.41%  '
$ws.Range('E49').Value = '  +5.30%  '
$ws.Range('D50').Value = '2.36'
$ws.Range('E50').Value = '  +8.44%  '
$ws.Range('D51').Value = "'7.00"
$ws.Range('E51').Value = '  +15.49%  '
